# "edited a text box"
#
# Slide 3 contains a Consolas-styled code/notes text box ("TextBox 38")
# listing the MigrationPath model and its functions. The edit removes the
# blank line that used to sit directly under the "MigrationPath" heading
# (merging it into the paragraph above) and, because the box has
# auto-fit-to-text (<a:spAutoFit/>) turned on, that shrinks the box and
# shifts it down a bit.

$p   = $ppt.ActivePresentation
$s   = $p.Slides.Item(3)
$shp = $s.Shapes.Item("TextBox 38")

$tf = $shp.TextFrame
$tr = $tf.TextRange

# Paragraph 1 = "MigrationPath", paragraph 2 = the now-removed blank line,
# paragraph 3 = "create_migration_path", etc. Deleting paragraph 2 joins
# the heading straight onto the following blank markup and (thanks to
# spAutoFit) automatically recalculates the shape's height.
$tr.Paragraphs(2).Delete()

# spAutoFit already fixed up the height; the top of the box still needs to
# move down to its new, slightly lower position (EMU 405332 -> 667572,
# i.e. points 31.9159 -> 52.56472).
$shp.Top = 52.564726409448824
